# Apply the "Append: 2026-01-12 06:33 JST" update to the Lancers worksheet.
#
# Summary of the change (per the target diff):
#  - Column B width 50 -> 47 (characters)
#  - Row 2: new scraped timestamp / title / URL / score / skills (row "moved up" one slot)
#  - Row 3: new scraped timestamp / title / price / URL / score / skills
#  - Row 4: only the timestamp refreshes, the rest of the row is unchanged
#  - Rows 5 and 6 (and their hyperlinks) are dropped entirely
#  - Sheet dimension shrinks from A1:H6 to A1:H4

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Column B width: 50 -> 47 characters -------------------------------
# This engine's ColumnWidth setter adds ~0.8333 (5/6) chars of internal
# padding before it gets persisted as the OOXML <col width>, so back that
# padding out up front to land on an exact stored width of 47.
$ws.Columns("B").ColumnWidth = 47 - (5/6)

# --- Drop all existing hyperlinks ---------------------------------------
# Hyperlinks.Delete() on this engine clears every hyperlink on the sheet
# regardless of which Range/Hyperlinks collection it was called on, so do
# it once up front and re-add only the ones that should survive.
$ws.Cells.Hyperlinks.Delete()

# --- Remove the two trailing scraped rows (5 and 6) ---------------------
$ws.Rows("5:6").Delete()

# --- Row 2: refreshed listing --------------------------------------------
$ws.Range("A2").Value = "2026-01-12 06:33:47"
$ws.Range("B2").Value = "大企業の業務効率化AIプロジェクトの技術方針策定を支援するAIテックリード募集"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5423720"
$ws.Range("G2").Value = 385
$ws.Range("H2").Value = "🔥AI,Ai ◆効率化"

# --- Row 3: refreshed listing --------------------------------------------
$ws.Range("A3").Value = "2026-01-12 06:33:47"
$ws.Range("B3").Value = "【Zapier設定のみ!作業時間~1時間】スプレッドシート・Gドライブ自動化構築(設計済)"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5469379"
$ws.Range("G3").Value = 255
$ws.Range("H3").Value = "🔥API ◆自動化"

# --- Row 4: only the "fetched at" timestamp changes ----------------------
$ws.Range("A4").Value = "2026-01-12 06:33:47"

# --- Re-create hyperlinks for the rows that remain (F2:F4) ---------------
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5423720")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5469379")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5469298")
